# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
# Swap the contents of columns B through AD between pairs of rows
# (row index A stays fixed; all other match data B..AD is exchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(105, 106),
    @(108, 110),
    @(113, 114),
    @(115, 117),
    @(118, 119),
    @(121, 122),
    @(139, 140)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($c = 2; $c -le 30; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value()
        $v2 = $ws.Cells.Item($r2, $c).Value()
        $ws.Cells.Item($r1, $c).Value = $v2
        $ws.Cells.Item($r2, $c).Value = $v1
    }
}
